{"js": "// The site footer block that used to precede the final page break was\n// removed on this rebuild: a blank paragraph, a page-break paragraph, and\n// the \"\u00a9 2020 . Contact: ...\" copyright paragraph are all deleted, leaving\n// the trailing blank + page-break paragraphs (already at the very end of\n// the document) directly after the \"LOB1036: Geometria Anal\u00edtica\n// (Requisito fraco)\" requirement line.\n\nconst body = context.document.body;\n\n// Locate the unique copyright/footer paragraph by its text.\nconst results = body.search(\"Contact: luizeleno@usp.br\", { matchCase: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find the copyright paragraph to remove.\");\n}\n\nconst copyrightPara = results.items[0].paragraphs.getFirst();\n\n// Walk backwards to the two empty paragraphs that immediately precede it\n// (a plain blank paragraph, and a blank page-break-before paragraph).\nconst blankPageBreak = copyrightPara.getPrevious();\nconst blankNormal = blankPageBreak.getPrevious();\n\n// Delete all three paragraphs (deepest/last one first so earlier\n// references stay valid).\ncopyrightPara.delete();\nblankPageBreak.delete();\nblankNormal.delete();\n\nawait context.sync();\n", "ps1": "# The site footer block that used to precede the final page break was\n# removed on this rebuild: a blank paragraph, a page-break paragraph, and\n# the \"\u00a9 2020 . Contact: ...\" copyright paragraph are all deleted, leaving\n# the trailing blank + page-break paragraphs (already at the very end of\n# the document) directly after the \"LOB1036: Geometria Analitica\n# (Requisito fraco)\" requirement line.\n\n$d = $word.ActiveDocument\n\n# Confirm the copyright/footer paragraph exists via Find.\n$findRange = $d.Content\n$findRange.Find.ClearFormatting()\n$findRange.Find.Text = \"Contact: luizeleno@usp.br\"\n$found = $findRange.Find.Execute()\n\nif (-not $found) {\n    Write-Output \"ERROR: copyright paragraph not found\"\n} else {\n    # Resolve the actual Paragraph object containing that text.\n    $target = $null\n    foreach ($p in $d.Paragraphs) {\n        if ($p.Range.Text -like \"*Contact: luizeleno@usp.br*\") {\n            $target = $p\n        }\n    }\n\n    # Walk backwards to the two empty paragraphs immediately preceding it\n    # (a blank page-break-before paragraph, then a plain blank paragraph).\n    $blankPageBreak = $target.Previous()\n    $blankNormal = $blankPageBreak.Previous()\n\n    # Delete all three paragraphs (the footer text, then the two blanks).\n    $target.Range.Delete()\n    $blankPageBreak.Range.Delete()\n    $blankNormal.Range.Delete()\n}\n"}
